$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For D/E (Price, Volume(1h)) cells we set NumberFormat to Text ("@") *before*
# assigning the new value so Excel stores the literal numeric-looking string
# (e.g. "332.07", "1.15%") instead of auto-converting it to a Number/Percentage.
# NumberFormat is set per-cell (not as a multi-area union range) since a union
# range only reliably applies the format to the first area.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "332.07"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "1.15%"

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "45.91"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "4.11%"

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.635"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.32%"

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.08367"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "4.42%"

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "2.060"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "4.17%"

$ws.Range("B7").Value = "MXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9848"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "3.82%"

$ws.Range("B8").Value = "BTSEToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "2.583"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-1.15%"

$ws.Range("B9").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C9").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1156"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "3.06%"

$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1926"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "3.58%"

$ws.Range("B11").Value = "MCDex"
$ws.Range("C11").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "10.41"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "-3.09%"

$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.09955"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "1.74%"

$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.04666"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "-0.22%"

$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.1059"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.61%"

$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001288"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.82%"

$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.006118"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "2.94%"

$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.375"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.52%"

$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "4.482"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "3.38%"

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "-3.14%"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1402"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.38%"

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.2654"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.40%"

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.04204"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "3.41%"

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001312"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "4.22%"

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.004629"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "6.79%"

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0001283"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.97%"

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0003747"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.10%"

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02768"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "7.35%"

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05807"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "2.97%"

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.007731"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "2.50%"

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.1435"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "2.70%"

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.007265"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-3.18%"

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002119"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "5.23%"

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.009054"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.01%"

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00007309"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "2.34%"

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.00000000751"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "0.24%"

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0005811"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "0.00%"

$ws.Range("B49").Value = "BOLO"
$ws.Range("C49").Value = "https://coinranking.com/coin/ogrGe0dEab+bolo-bolo"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.003498"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-2.72%"

$ws.Range("B50").Value = "CoinbaseStockToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/_ZA6fIr53+coinbasestocktoken-coin"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.003505"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.68%"

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.00002104"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "0.24%"
